# Rename "task book" / "prevTaskBook" references to "organizer" / "prevOrganizer"
# across the UndoRedoExecuteUndoStackDiagram slide (Addressbook -> PrioriTask rename).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The four small "previous state" tables (Table 45, Table 43, Table 42, Table 44)
# each contain a cell whose second paragraph reads "prevTaskBook = s2" (or "= s3").
# Rename the "prevTaskBook" token to "prevOrganizer" in each, leaving the rest of
# the cell text untouched.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable -eq -1) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cellShape = $tbl.Cell($r, $c).Shape
                $tr = $cellShape.TextFrame.TextRange
                $full = $tr.Text
                $idx = $full.IndexOf("prevTaskBook")
                if ($idx -ge 0) {
                    $sub = $tr.Characters($idx + 1, 12)
                    $sub.Text = "prevOrganizer"
                }
            }
        }
    }
}

# The caption textbox: "The state of the task book (before 'add ... " ->
# "The state of the organizer (before 'add ... "
$caption = $s.Shapes.Item("TextBox 1")
$capRange = $caption.TextFrame.TextRange
$capText = $capRange.Text
$capIdx = $capText.IndexOf("task book ")
if ($capIdx -ge 0) {
    $capSub = $capRange.Characters($capIdx + 1, 10)
    $capSub.Text = "organizer "
}
